# Generate Report for Handoff
#
# The localization file "9b4c2133-3d5d-4c3a-885a-eec3504d48c5.md" has
# completed its handoff step, so every sheet's row for that file is
# updated to reflect the new "Ready for handoff" status, the downgraded
# ("mt") priority and the refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
# Row 3 is the 9b4c2133-...md file (row 2 is 077c6bd7-...md, untouched)
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-12 02:36:19"
$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
# Row 3 is the 9b4c2133-...md file (row 2 is 077c6bd7-...md, untouched)
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-12 02:36:13"
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
# Row 3 is the 9b4c2133-...md file (row 2 is 077c6bd7-...md, untouched)
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-12 02:36:19"
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
